$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.767.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.475.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.28%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +16.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.467.21"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.689"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +28.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.54"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.98%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.017.06"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.03"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.503.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.612.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.54%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.88"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +22.33%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "309.52"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.181"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.48%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.88"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "43.12"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.64%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0492"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.62"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.94%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.66"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.60%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.125"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.40"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.285"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.33"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.205.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.810.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.52%  "
